$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.203.34"

$ws.Range("D3").Value = "2.025.40"
$ws.Range("E3").Value = "  +3.20%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'247.68"
$ws.Range("D5").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "'0.630"
$ws.Range("D6").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E6").Value = "  +3.05%  "

$ws.Range("D7").Value = "'60.39"
$ws.Range("D7").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E7").Value = "  -2.24%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.394"
$ws.Range("D9").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E9").Value = "  +4.78%  "

$ws.Range("D10").Value = "'0.0813"
$ws.Range("D10").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E10").Value = "  +2.11%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "'15.28"
$ws.Range("D12").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E12").Value = "  +6.93%  "

$ws.Range("D13").Value = "'22.55"
$ws.Range("D13").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E13").Value = "  +1.64%  "

$ws.Range("E14").Value = "  +3.10%  "

$ws.Range("D15").Value = "2.322.18"
$ws.Range("E15").Value = "  +3.24%  "

$ws.Range("E16").Value = "  +3.99%  "

$ws.Range("D17").Value = "2.027.44"
$ws.Range("E17").Value = "  +2.91%  "

$ws.Range("D18").Value = "37.157.35"
$ws.Range("E18").Value = "  +1.46%  "

$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("D20").Value = "0.0₃0869"

$ws.Range("E21").Value = "  +3.57%  "

$ws.Range("D22").Value = "'231.27"
$ws.Range("D22").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  +2.79%  "

$ws.Range("D27").Value = "'163.97"
$ws.Range("D27").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E27").Value = "  +2.09%  "

$ws.Range("E28").Value = "  -3.79%  "

$ws.Range("D29").Value = "'19.86"
$ws.Range("D29").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E29").Value = "  +2.16%  "

$ws.Range("D30").Value = "'1.39"
$ws.Range("D30").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E30").Value = "  +6.94%  "

$ws.Range("E31").Value = "  +1.95%  "

$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("D33").Value = "'0.0665"
$ws.Range("D33").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E33").Value = "  +8.03%  "

$ws.Range("D34").Value = "'4.56"
$ws.Range("D34").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("D35").Value = "'2.49"
$ws.Range("D35").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E35").Value = "  +9.49%  "

$ws.Range("D36").Value = "'3.45"
$ws.Range("D36").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E36").Value = "  -3.68%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").Value = "'1.81"
$ws.Range("D38").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("E39").Value = "  -1.71%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("E43").Value = "  +1.97%  "

$ws.Range("D44").Value = "'16.83"
$ws.Range("D44").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E44").Value = "  +4.61%  "

$ws.Range("D45").Value = "'92.26"
$ws.Range("D45").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E45").Value = "  +3.80%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "'1.07"
$ws.Range("D46").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E46").Value = "  +3.36%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.391.24"
$ws.Range("E47").Value = "  +1.42%  "

$ws.Range("D48").Value = "'7.50"
$ws.Range("D48").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E48").Value = "  +4.88%  "

$ws.Range("D49").Value = "'2.19"
$ws.Range("D49").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E49").Value = "  +18.61%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").Value = "'46.90"
$ws.Range("D51").Style = "Normal"  # keep as plain text, no quote-prefix style
$ws.Range("E51").Value = "  +2.97%  "
